$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 8 (ano=2025) with refreshed Bibi data
$ws.Range("C8").Value = 874
$ws.Range("E8").Value = 728
$ws.Range("G8").Value = 83.29519450800915
$ws.Range("H8").Value = 16.70480549199085
